# Weekly driver report update for 2025-04-21
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Driver Summary")

# Update "Bad Drivers" table values per the weekly refresh
$ws.Range("D3").Value = 96.59999999999999
$ws.Range("C4").Value = 352
$ws.Range("D5").Value = 97.7
$ws.Range("C6").Value = 415
